$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-05-12 Sunday" "2024-05-13 Monday"

Replace-Text "25×38=" "89×28="
Replace-Text "32×50=" "80×62="
Replace-Text "25×73=" "67×80="
Replace-Text "63×29=" "44×92="
Replace-Text "29×52=" "87×82="
Replace-Text "51×98=" "85×50="
Replace-Text "57×68=" "14×79="
Replace-Text "75×12=" "60×70="
Replace-Text "47×38=" "20×51="
Replace-Text "94×84=" "27×78="
Replace-Text "76×64=" "57×66="
Replace-Text "23×58=" "61×45="
Replace-Text "11×69=" "67×30="
Replace-Text "49×27=" "63×16="
Replace-Text "81×33=" "34×51="
Replace-Text "63×83=" "17×17="
Replace-Text "86×72=" "76×85="
Replace-Text "14×31=" "58×50="
Replace-Text "29×96=" "73×73="
Replace-Text "88×31=" "80×12="
Replace-Text "83×24=" "37×30="
Replace-Text "33×47=" "12×82="
Replace-Text "38×48=" "42×83="
Replace-Text "60×17=" "81×16="
Replace-Text "37×37=" "14×27="
